$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

$ws.Range("H43").Value = 1379.8
$ws.Range("I43").Value = 849.5
$ws.Range("J43").Value = 1733.3334
$ws.Range("K43").Value = 849.5
$ws.Range("L43").Value = 1733.3334
$ws.Range("M43").Value = -780.5

$ws.Range("H64").Value = 7483.3335
$ws.Range("I64").Value = 5450
$ws.Range("J64").Value = 8500
$ws.Range("K64").Value = 5450
$ws.Range("L64").Value = 8500
$ws.Range("M64").Value = -5202
$ws.Range("N64").Value = -8996

$ws.Range("H67").Value = 7483.3335
$ws.Range("I67").Value = 5450
$ws.Range("J67").Value = 8500
$ws.Range("K67").Value = 5450
$ws.Range("L67").Value = 8500
$ws.Range("M67").Value = -4592
$ws.Range("N67").Value = -10216

$ws.Range("H107").Value = 549.75
$ws.Range("I107").Value = 549.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 549.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1370.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1500
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1387
$ws.Range("N2").ClearContents()

$ws.Range("H116").Value = 1500
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 794
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1500
$ws.Range("I3").Value = 1500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1386
$ws.Range("N3").ClearContents()

$ws.Range("H20").Value = 4242.2856
$ws.Range("I20").Value = 4405.6
$ws.Range("J20").Value = 3834
$ws.Range("K20").Value = 4405.6
$ws.Range("L20").Value = 3834
$ws.Range("M20").Value = -4158.6
$ws.Range("N20").Value = -4328

$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 100
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 73
$ws.Range("N22").Value = -846

$ws.Range("H29").Value = 1185.1428
$ws.Range("I29").Value = 1369.2
$ws.Range("J29").Value = 725
$ws.Range("K29").Value = 1369.2
$ws.Range("L29").Value = 725
$ws.Range("M29").Value = -1080.2
$ws.Range("N29").Value = -1303

$ws.Range("H36").Value = 991.8461
$ws.Range("I36").Value = 815.36365
$ws.Range("J36").Value = 1962.5
$ws.Range("K36").Value = 815.36365
$ws.Range("L36").Value = 1962.5
$ws.Range("M36").Value = -281.36365

$ws.Range("H99").Value = 142858060
$ws.Range("I99").Value = 166667680
$ws.Range("J99").Value = 400
$ws.Range("K99").Value = 166667680
$ws.Range("L99").Value = 400
$ws.Range("M99").Value = -166666182
$ws.Range("N99").Value = -3396

$ws.Range("H105").Value = 7577087
$ws.Range("I105").Value = 15152708
$ws.Range("J105").Value = 1466.6666
$ws.Range("K105").Value = 15152708
$ws.Range("L105").Value = 1466.6666
$ws.Range("M105").Value = -15150961

$ws.Range("H107").Value = 41672056
$ws.Range("I107").Value = 83335370
$ws.Range("J107").Value = 8750
$ws.Range("K107").Value = 83335370
$ws.Range("L107").Value = 8750
$ws.Range("M107").Value = -83333450

$ws.Range("H134").Value = 3044.9285
$ws.Range("I134").Value = 923.3333
$ws.Range("J134").Value = 15774.5
$ws.Range("K134").Value = 2769.9999
$ws.Range("L134").Value = 47323.5
$ws.Range("M134").Value = -234.9998999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 547.5
$ws.Range("I2").Value = 577.5
$ws.Range("J2").Value = 457.5
$ws.Range("K2").Value = 577.5
$ws.Range("L2").Value = 457.5
$ws.Range("M2").Value = -464.5
$ws.Range("N2").Value = -683.5

$ws.Range("H25").Value = 2525
$ws.Range("I25").Value = 2830
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 2830
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = -2656
$ws.Range("N25").Value = -1348

$ws.Range("H31").Value = 6563.515
$ws.Range("I31").Value = 4132.6665
$ws.Range("J31").Value = 7103.7036
$ws.Range("K31").Value = 4132.6665
$ws.Range("L31").Value = 7103.7036
$ws.Range("M31").Value = -3837.6665

$ws.Range("H34").Value = 6563.515
$ws.Range("I34").Value = 4132.6665
$ws.Range("J34").Value = 7103.7036
$ws.Range("K34").Value = 4132.6665
$ws.Range("L34").Value = 7103.7036
$ws.Range("M34").Value = -3930.6665

$ws.Range("H51").Value = 27789

$ws.Range("H61").Value = 27789

$ws.Range("H134").Value = 2661.647
$ws.Range("I134").Value = 1092.7273
$ws.Range("J134").Value = 5538
$ws.Range("K134").Value = 3278.1819
$ws.Range("L134").Value = 16614
$ws.Range("M134").Value = -743.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 277.5
$ws.Range("I12").Value = 232.4
$ws.Range("J12").Value = 298
$ws.Range("K12").Value = 697.2
$ws.Range("L12").Value = 894
$ws.Range("M12").Value = -524.2

$ws.Range("H39").Value = 7781
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 7781
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 23343
$ws.Range("N39").Value = -23931

$ws.Range("H92").Value = 3999.7144
$ws.Range("I92").Value = 1266.3334
$ws.Range("J92").Value = 6049.75
$ws.Range("K92").Value = 3799.0002
$ws.Range("L92").Value = 18149.25
$ws.Range("M92").Value = -2551.0002

$ws.Range("H97").Value = 2394.3333
$ws.Range("I97").Value = 381.5
$ws.Range("J97").Value = 6420
$ws.Range("K97").Value = 1144.5
$ws.Range("L97").Value = 19260
$ws.Range("M97").Value = -648.5
$ws.Range("N97").Value = -20252

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 30999.8
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 30999.8
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 30999.8
$ws.Range("N26").Value = -31559.8

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H50").Value = 30999.8
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 30999.8
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 30999.8
$ws.Range("N50").Value = -31995.8

$ws.Range("H58").Value = 27500
$ws.Range("I58").Value = 30000
$ws.Range("J58").Value = 25000
$ws.Range("K58").Value = 30000
$ws.Range("L58").Value = 25000
$ws.Range("M58").Value = -29723
$ws.Range("N58").Value = -25554

$ws.Range("H80").Value = 3159.6
$ws.Range("I80").Value = 3159.6
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3159.6
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2161.6
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 3159.6
$ws.Range("I83").Value = 3159.6
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15798
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10806
$ws.Range("N83").ClearContents()

$ws.Range("H102").Value = 1259.7894
$ws.Range("I102").Value = 1172.7646
$ws.Range("J102").Value = 1999.5
$ws.Range("K102").Value = 1172.7646
$ws.Range("L102").Value = 1999.5
$ws.Range("M102").Value = 449.2354
$ws.Range("N102").Value = -5243.5

$ws.Range("H122").Value = 2325.2
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 2656.5
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 7969.5
$ws.Range("M122").Value = -550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5497
$ws.Range("I7").Value = 3996.3333
$ws.Range("J7").Value = 9999
$ws.Range("K7").Value = 3996.3333
$ws.Range("L7").Value = 9999
$ws.Range("M7").Value = -3884.3333
$ws.Range("N7").Value = -10223

$ws.Range("H16").Value = 1178.4
$ws.Range("I16").Value = 1296.6666
$ws.Range("J16").Value = 1001
$ws.Range("K16").Value = 1296.6666
$ws.Range("L16").Value = 1001
$ws.Range("M16").Value = -1126.6666

$ws.Range("H22").Value = 1819.1111
$ws.Range("I22").Value = 1633.3334
$ws.Range("J22").Value = 1912
$ws.Range("K22").Value = 1633.3334
$ws.Range("L22").Value = 1912
$ws.Range("M22").Value = -1338.3334
$ws.Range("N22").Value = -2502

$ws.Range("H27").Value = 1819.1111
$ws.Range("I27").Value = 1633.3334
$ws.Range("J27").Value = 1912
$ws.Range("K27").Value = 1633.3334
$ws.Range("L27").Value = 1912
$ws.Range("M27").Value = -1526.3334
$ws.Range("N27").Value = -2126

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 809.75
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 846.3333
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 846.3333
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -1222.3333

$ws.Range("H64").Value = 35000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 35000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35450

$ws.Range("H67").Value = 35000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 35000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36560

$ws.Range("H68").Value = 6788
$ws.Range("I68").Value = 4856.7144
$ws.Range("J68").Value = 9041.166999999999
$ws.Range("K68").Value = 4856.7144
$ws.Range("L68").Value = 9041.166999999999
$ws.Range("M68").Value = -4107.7144
$ws.Range("N68").Value = -10539.167

$ws.Range("H71").Value = 6788
$ws.Range("I71").Value = 4856.7144
$ws.Range("J71").Value = 9041.166999999999
$ws.Range("K71").Value = 24283.572
$ws.Range("L71").Value = 45205.835
$ws.Range("M71").Value = -20539.572
$ws.Range("N71").Value = -52693.835

$ws.Range("H122").Value = 999.5
$ws.Range("I122").Value = 999.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2998.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -548.5

$ws.Range("H126").Value = 5497
$ws.Range("I126").Value = 3996.3333
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 11988.9999
$ws.Range("L126").Value = 29997
$ws.Range("M126").Value = -9518.999899999999
$ws.Range("N126").Value = -34937

$ws.Range("H136").Value = 2857.8333
$ws.Range("I136").Value = 2648.9
$ws.Range("J136").Value = 3902.5
$ws.Range("K136").Value = 7946.700000000001
$ws.Range("L136").Value = 11707.5
$ws.Range("M136").Value = -5396.700000000001

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 10199.4
$ws.Range("I14").Value = 15000
$ws.Range("J14").Value = 8999.25
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 8999.25
$ws.Range("M14").Value = -14832
$ws.Range("N14").Value = -9335.25

$ws.Range("H122").Value = 3569.3076
$ws.Range("I122").Value = 2759.5715
$ws.Range("J122").Value = 4514
$ws.Range("K122").Value = 8278.7145
$ws.Range("L122").Value = 13542
$ws.Range("M122").Value = -5828.7145
